$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update player names (column B) - same row positions, new names
$ws.Range("B2").Value = "Hîstorîa"
$ws.Range("B3").Value = "mommy cat"
$ws.Range("B4").Value = "Zeddxx"
$ws.Range("B5").Value = "xiiixo"
$ws.Range("B6").Value = "XENOX"
$ws.Range("B7").Value = "Faizenissobad"
$ws.Range("B8").Value = "DrLimits"
$ws.Range("B9").Value = "Lazar"
$ws.Range("B10").Value = "Bones Slayer Zed"
$ws.Range("B11").Value = "fenix3006"

# Update agent (column C)
$ws.Range("C2").Value = "Raze"
$ws.Range("C3").Value = "Sage"
$ws.Range("C4").Value = "Reyna"
$ws.Range("C5").Value = "Sage"
$ws.Range("C6").Value = "Sova"
$ws.Range("C7").Value = "Chamber"
$ws.Range("C8").Value = "Neon"
$ws.Range("C9").Value = "Brimstone"
$ws.Range("C10").Value = "Jett"
$ws.Range("C11").Value = "Jett"

# Update stats columns:
# D=Scores, E=Total Damage, F=Total Damage received, G=Headshots,
# H=Kills, I=Deaths, J=Ultimate, K=KD

$ws.Range("D2").Value = 2336
$ws.Range("E2").Value = 1739
$ws.Range("F2").Value = 2847
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 0.5714285714285714

$ws.Range("D3").Value = 4716
$ws.Range("E3").Value = 2836
$ws.Range("F3").Value = 3269
$ws.Range("G3").Value = 9
$ws.Range("H3").Value = 16
$ws.Range("I3").Value = 17
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.9411764705882353

$ws.Range("D4").Value = 6634
$ws.Range("E4").Value = 3749
$ws.Range("F4").Value = 2258
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 25
$ws.Range("I4").Value = 12
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2.083333333333333

$ws.Range("D5").Value = 4085
$ws.Range("E5").Value = 2406
$ws.Range("F5").Value = 2645
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 14
$ws.Range("I5").Value = 15
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0.9333333333333333

$ws.Range("D6").Value = 4540
$ws.Range("E6").Value = 2820
$ws.Range("F6").Value = 2768
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 15
$ws.Range("I6").Value = 16
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 0.9375

$ws.Range("D7").Value = 5140
$ws.Range("E7").Value = 3465
$ws.Range("F7").Value = 2989
$ws.Range("G7").Value = 15
$ws.Range("H7").Value = 19
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = 4
$ws.Range("K7").Value = 1.266666666666667

$ws.Range("D8").Value = 2243
$ws.Range("E8").Value = 1660
$ws.Range("F8").Value = 2780
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 17
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 0.3529411764705883

$ws.Range("D9").Value = 3933
$ws.Range("E9").Value = 2743
$ws.Range("F9").Value = 2694
$ws.Range("G9").Value = 11
$ws.Range("H9").Value = 15
$ws.Range("I9").Value = 15
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 1

$ws.Range("D10").Value = 4038
$ws.Range("E10").Value = 2444
$ws.Range("F10").Value = 2397
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 15
$ws.Range("I10").Value = 12
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = 1.25

$ws.Range("D11").Value = 4969
$ws.Range("E11").Value = 3010
$ws.Range("F11").Value = 3067
$ws.Range("G11").Value = 12
$ws.Range("H11").Value = 17
$ws.Range("I11").Value = 17
$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 1
